# Auto-generated edit script: updates crypto price/volume table (cryptos.xlsx)
# Commit: "Updated cryptos list on Sat Apr 15 04:27:50 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.455.53'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '2.091.74'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.23%  '
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5206'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4364'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.96'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +16.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08831'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.17'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.70%  '
$ws.Range('D13').Value = '2.089.40'
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.672'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.656'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '95.68'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.79%  '
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001118'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06592'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.255'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.46%  '
$ws.Range('D23').Value = '30.510.22'
$ws.Range('E23').Value = '  -1.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.341'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = '2.333.73'
$ws.Range('E26').Value = '  -1.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.543'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '131.36'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.180'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.643'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.136'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.905'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02569'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06796'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.437'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.27%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.59'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2246'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6855'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.259'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.94'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6326'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('E47').Value = '  -2.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.627'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.15%  '
$ws.Range('B49').Value = 'WEMIXTOKEN'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.234'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.44%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.239'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.41'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.31%  '
